# machineindex -> machineno 기준으로 변경 및 machineindex 컬럼 삭제.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "machine_master" to "Sheet1"
$ws.Name = "Sheet1"

# Delete column A (machineindex); machineno/machinename shift left into A/B
$ws.Columns.Item(1).Delete()
